# Updated cryptos list on Sat Oct 21 17:24:05 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns with newly scraped
# values, and fixes the Aave / BitcoinSV row ordering (rows 44-45).
#
# Numeric-looking price strings are written with a leading apostrophe so
# Excel stores them as literal text (matching the source workbook's
# inline-string cells) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.650.86"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.613.41"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "`'212.55"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "`'0.993"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "`'29.09"
$ws.Range("E8").Value = "  +9.60%  "
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "1.843.90"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "1.624.98"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "`'0.569"
$ws.Range("E14").Value = "  +6.88%  "
$ws.Range("D15").Value = "`'3.88"
$ws.Range("E15").Value = "  +5.33%  "
$ws.Range("D16").Value = "29.664.49"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "`'8.99"
$ws.Range("E17").Value = "  +17.72%  "
$ws.Range("D18").Value = "`'64.15"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "`'240.71"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "0.0₃0710"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "`'4.10"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "`'9.64"
$ws.Range("E23").Value = "  +5.76%  "
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").Value = "`'156.50"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "`'15.62"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").Value = "`'6.59"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  +3.40%  "
$ws.Range("D31").Value = "`'1.09"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "`'3.21"
$ws.Range("E33").Value = "  +3.79%  "
$ws.Range("D34").Value = "1.425.84"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Value = "`'2.86"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "`'2.28"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D40").Value = "`'0.556"
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("D43").Value = "`'1.96"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "`'53.94"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "`'69.72"
$ws.Range("E45").Value = "  +6.37%  "
$ws.Range("D46").Value = "`'0.993"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  +17.79%  "
$ws.Range("D48").Value = "`'5.44"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").Value = "1.753.42"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").Value = "`'87.81"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "0.0₆0109"
$ws.Range("E51").Value = "  +7.95%  "
